# LMS-2523 Update BaSynthec Validation
# Rename the MGP253/MGP776 strain labels on the "openbis-data" sheet to
# their new "JJS-" prefixed identifiers, and leave the sheet selection on
# "openbis-data" (as the author did while making this edit).

$wb = $excel.ActiveWorkbook
$wsData = $wb.Worksheets.Item("openbis-data")

$wsData.Range("D1").Value = "JJS-MGP253"
$wsData.Range("E1").Value = "JJS-MGP776"

$wsData.Activate()
$wsData.Range("I16").Select()
